$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

try {
  $excel.Goto($ws.Range("A128"), $true)
  Write-Host "Goto worked"
} catch {
  Write-Host "Goto failed: $_"
}
$ws.Range("A128:XFD128").Select()
